$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Comments
# ---------------------------------------------------------------------------
# Comment 0: "The program can run on recent versions of Windows, Linux, and Mac"
$p1 = $d.Paragraphs.Item(2)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$d.Comments.Add($r1, "Only Windows and Linux are available for now." + [char]13 + [char]13 + "Known bug: When resizing the window on the Windows release, scroll bars appear around the scene.")

# Comment 1: "A select tool is provided in the form of a lasso select..."
$p2 = $d.Paragraphs.Item(17)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$d.Comments.Add($r2, "Known bug: When a domain is deleted which is merged with another domain, the other domain" + [char]8217 + "s walls are not updated.")

# Comment 2: "A domain can be resized by dragging one of its walls"
$p3 = $d.Paragraphs.Item(24)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$d.Comments.Add($r3, "Known bug: When dragging a wall across its opposite wall, they cannot be separated anymore.")

# ---------------------------------------------------------------------------
# 2) Text corrections
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("rectangles (or lines) on the GUI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "rectangles on the GUI", 2)

$d.Content.Find.Execute("A domain or wall can be deleted from the scene", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A domain can be deleted from the scene", 2)

# ---------------------------------------------------------------------------
# 3) Colour re-coding (orange FF7F00 -> green 00A800) on several bullets
# ---------------------------------------------------------------------------
# "The scene can be moved by dragging the window background" -- whole paragraph (incl. pilcrow) turns green
$d.Paragraphs.Item(4).Range.Font.Color = 0x00A800

# "The scene can be zoomed by scrolling" -- only the run text turns green, pilcrow stays orange
$p = $d.Paragraphs.Item(5)
$d.Range($p.Range.Start, $p.Range.End - 1).Font.Color = 0x00A800

# "Dimensions of the walls are shown as a number next to the wall (on the outside where possible)"
$p = $d.Paragraphs.Item(23)
$d.Range($p.Range.Start, $p.Range.End - 1).Font.Color = 0x00A800

# "Two touching domains can have their common wall removed"
$p = $d.Paragraphs.Item(26)
$d.Range($p.Range.Start, $p.Range.End - 1).Font.Color = 0x00A800

# "A domain can be deleted from the scene" (renamed paragraph) -- whole paragraph turns green
$d.Paragraphs.Item(27).Range.Font.Color = 0x00A800

# "A wall has an absorption coefficient that can be altered" -- whole paragraph gains green colour
$d.Paragraphs.Item(28).Range.Font.Color = 0x00A800

# "Drawing or resizing a domain anchors to the grid"
$p = $d.Paragraphs.Item(29)
$d.Range($p.Range.Start, $p.Range.End - 1).Font.Color = 0x00A800

# "A source or receiver can be relocated by dragging it around" -- whole paragraph turns green
$d.Paragraphs.Item(32).Range.Font.Color = 0x00A800

# "A source or receiver can be deleted from the scene"
$p = $d.Paragraphs.Item(33)
$d.Range($p.Range.Start, $p.Range.End - 1).Font.Color = 0x00A800

# "Drawing or relocating a source or receiver anchors to the grid"
$p = $d.Paragraphs.Item(34)
$d.Range($p.Range.Start, $p.Range.End - 1).Font.Color = 0x00A800
